# Update the "I will be collaborating with the following company/department"
# answer on the DataEntry sheet from "N/A" to "Conference Services ".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataEntry")

# Update the collaborator cell (B8) value.
$ws.Range("B8").Value = "Conference Services "

# Update the active selection to match the edited cell.
$ws.Activate()
$ws.Range("B8").Select()

# Recalculate so the DataBase sheet formula referencing DataEntry!B8 updates.
$excel.Calculate()
